# Update the LR-pair TPM table (rows 2-10) with the recalculated values.
# The original sheet had 6 data rows (ECs/FAPs pairs only); the refreshed
# TPM run adds the "MuSCs" cluster, producing a full 3x3 cluster matrix
# (9 data rows, A1:T10) with updated expression / specificity figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Amelx"
$ws.Cells.Item(2,3).Value = "Lamp1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1.0
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.052706
$ws.Cells.Item(2,8).Value = 0.158118
$ws.Cells.Item(2,9).Value = 0.0131977395622021
$ws.Cells.Item(2,10).Value = 0.0131977395622021
$ws.Cells.Item(2,11).Value = 3.0
$ws.Cells.Item(2,12).Value = 1.0
$ws.Cells.Item(2,13).Value = 32.130737
$ws.Cells.Item(2,14).Value = 96.392211
$ws.Cells.Item(2,15).Value = 0.1572847769351211
$ws.Cells.Item(2,16).Value = 0.1572847769351211
$ws.Cells.Item(2,17).Value = 1.693482624322
$ws.Cells.Item(2,18).Value = 15.241343618898
$ws.Cells.Item(2,19).Value = 0.00207580352308878
$ws.Cells.Item(2,20).Value = 0.00207580352308878

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Amelx"
$ws.Cells.Item(3,3).Value = "Lamp1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1.0
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.052706
$ws.Cells.Item(3,8).Value = 0.158118
$ws.Cells.Item(3,9).Value = 0.0131977395622021
$ws.Cells.Item(3,10).Value = 0.0131977395622021
$ws.Cells.Item(3,11).Value = 3.0
$ws.Cells.Item(3,12).Value = 1.0
$ws.Cells.Item(3,13).Value = 124.5871326666667
$ws.Cells.Item(3,14).Value = 373.761398
$ws.Cells.Item(3,15).Value = 0.6098727013470931
$ws.Cells.Item(3,16).Value = 0.6098727013470933
$ws.Cells.Item(3,17).Value = 6.566489414329333
$ws.Cells.Item(3,18).Value = 59.098404728964
$ws.Cells.Item(3,19).Value = 0.0080489410784756
$ws.Cells.Item(3,20).Value = 0.008048941078475601

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Amelx"
$ws.Cells.Item(4,3).Value = "Lamp1"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 1.0
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.052706
$ws.Cells.Item(4,8).Value = 0.158118
$ws.Cells.Item(4,9).Value = 0.0131977395622021
$ws.Cells.Item(4,10).Value = 0.0131977395622021
$ws.Cells.Item(4,11).Value = 3.0
$ws.Cells.Item(4,12).Value = 1.0
$ws.Cells.Item(4,13).Value = 47.56596266666667
$ws.Cells.Item(4,14).Value = 142.697888
$ws.Cells.Item(4,15).Value = 0.2328425217177857
$ws.Cells.Item(4,16).Value = 0.2328425217177857
$ws.Cells.Item(4,17).Value = 2.507011628309334
$ws.Cells.Item(4,18).Value = 22.563104654784
$ws.Cells.Item(4,19).Value = 0.003072994960637724
$ws.Cells.Item(4,20).Value = 0.003072994960637724

# Row 5: FAPs -> ECs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Amelx"
$ws.Cells.Item(5,3).Value = "Lamp1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3.0
$ws.Cells.Item(5,6).Value = 1.0
$ws.Cells.Item(5,7).Value = 1.388571333333333
$ws.Cells.Item(5,8).Value = 4.165713999999999
$ws.Cells.Item(5,9).Value = 0.347702402399595
$ws.Cells.Item(5,10).Value = 0.347702402399595
$ws.Cells.Item(5,11).Value = 3.0
$ws.Cells.Item(5,12).Value = 1.0
$ws.Cells.Item(5,13).Value = 32.130737
$ws.Cells.Item(5,14).Value = 96.392211
$ws.Cells.Item(5,15).Value = 0.1572847769351211
$ws.Cells.Item(5,16).Value = 0.1572847769351211
$ws.Cells.Item(5,17).Value = 44.61582031707267
$ws.Cells.Item(5,18).Value = 401.542382853654
$ws.Cells.Item(5,19).Value = 0.05468829480122601
$ws.Cells.Item(5,20).Value = 0.05468829480122601

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Amelx"
$ws.Cells.Item(6,3).Value = "Lamp1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3.0
$ws.Cells.Item(6,6).Value = 1.0
$ws.Cells.Item(6,7).Value = 1.388571333333333
$ws.Cells.Item(6,8).Value = 4.165713999999999
$ws.Cells.Item(6,9).Value = 0.347702402399595
$ws.Cells.Item(6,10).Value = 0.347702402399595
$ws.Cells.Item(6,11).Value = 3.0
$ws.Cells.Item(6,12).Value = 1.0
$ws.Cells.Item(6,13).Value = 124.5871326666667
$ws.Cells.Item(6,14).Value = 373.761398
$ws.Cells.Item(6,15).Value = 0.6098727013470931
$ws.Cells.Item(6,16).Value = 0.6098727013470933
$ws.Cells.Item(6,17).Value = 172.9981209231302
$ws.Cells.Item(6,18).Value = 1556.983088308172
$ws.Cells.Item(6,19).Value = 0.212054203416315
$ws.Cells.Item(6,20).Value = 0.2120542034163151

# Row 7: FAPs -> MuSCs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Amelx"
$ws.Cells.Item(7,3).Value = "Lamp1"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 3.0
$ws.Cells.Item(7,6).Value = 1.0
$ws.Cells.Item(7,7).Value = 1.388571333333333
$ws.Cells.Item(7,8).Value = 4.165713999999999
$ws.Cells.Item(7,9).Value = 0.347702402399595
$ws.Cells.Item(7,10).Value = 0.347702402399595
$ws.Cells.Item(7,11).Value = 3.0
$ws.Cells.Item(7,12).Value = 1.0
$ws.Cells.Item(7,13).Value = 47.56596266666667
$ws.Cells.Item(7,14).Value = 142.697888
$ws.Cells.Item(7,15).Value = 0.2328425217177857
$ws.Cells.Item(7,16).Value = 0.2328425217177857
$ws.Cells.Item(7,17).Value = 66.04873220133689
$ws.Cells.Item(7,18).Value = 594.438589812032
$ws.Cells.Item(7,19).Value = 0.08095990418205398
$ws.Cells.Item(7,20).Value = 0.08095990418205398

# Row 8: MuSCs -> ECs
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Amelx"
$ws.Cells.Item(8,3).Value = "Lamp1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3.0
$ws.Cells.Item(8,6).Value = 1.0
$ws.Cells.Item(8,7).Value = 2.552285333333333
$ws.Cells.Item(8,8).Value = 7.656856
$ws.Cells.Item(8,9).Value = 0.6390998580382028
$ws.Cells.Item(8,10).Value = 0.6390998580382028
$ws.Cells.Item(8,11).Value = 3.0
$ws.Cells.Item(8,12).Value = 1.0
$ws.Cells.Item(8,13).Value = 32.130737
$ws.Cells.Item(8,14).Value = 96.392211
$ws.Cells.Item(8,15).Value = 0.1572847769351211
$ws.Cells.Item(8,16).Value = 0.1572847769351211
$ws.Cells.Item(8,17).Value = 82.00680879429068
$ws.Cells.Item(8,18).Value = 738.0612791486161
$ws.Cells.Item(8,19).Value = 0.1005206786108063
$ws.Cells.Item(8,20).Value = 0.1005206786108063

# Row 9: MuSCs -> FAPs
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Amelx"
$ws.Cells.Item(9,3).Value = "Lamp1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3.0
$ws.Cells.Item(9,6).Value = 1.0
$ws.Cells.Item(9,7).Value = 2.552285333333333
$ws.Cells.Item(9,8).Value = 7.656856
$ws.Cells.Item(9,9).Value = 0.6390998580382028
$ws.Cells.Item(9,10).Value = 0.6390998580382028
$ws.Cells.Item(9,11).Value = 3.0
$ws.Cells.Item(9,12).Value = 1.0
$ws.Cells.Item(9,13).Value = 124.5871326666667
$ws.Cells.Item(9,14).Value = 373.761398
$ws.Cells.Item(9,15).Value = 0.6098727013470931
$ws.Cells.Item(9,16).Value = 0.6098727013470933
$ws.Cells.Item(9,17).Value = 317.9819114271875
$ws.Cells.Item(9,18).Value = 2861.837202844688
$ws.Cells.Item(9,19).Value = 0.3897695568523025
$ws.Cells.Item(9,20).Value = 0.3897695568523026

# Row 10: MuSCs -> MuSCs
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Amelx"
$ws.Cells.Item(10,3).Value = "Lamp1"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 3.0
$ws.Cells.Item(10,6).Value = 1.0
$ws.Cells.Item(10,7).Value = 2.552285333333333
$ws.Cells.Item(10,8).Value = 7.656856
$ws.Cells.Item(10,9).Value = 0.6390998580382028
$ws.Cells.Item(10,10).Value = 0.6390998580382028
$ws.Cells.Item(10,11).Value = 3.0
$ws.Cells.Item(10,12).Value = 1.0
$ws.Cells.Item(10,13).Value = 47.56596266666667
$ws.Cells.Item(10,14).Value = 142.697888
$ws.Cells.Item(10,15).Value = 0.2328425217177857
$ws.Cells.Item(10,16).Value = 0.2328425217177857
$ws.Cells.Item(10,17).Value = 121.4019088800142
$ws.Cells.Item(10,18).Value = 1092.617179920128
$ws.Cells.Item(10,19).Value = 0.148809622575094
$ws.Cells.Item(10,20).Value = 0.148809622575094
